# Project alert1 is saved.TEST Author: admin. Type: SAVE.
# The only functional change in this revision is the numeric value of
# cell C10 on the "Rules" sheet, which moves from 18 to 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 19
